$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 12503757
$ws.Range("I86").Value = 1733.6666
$ws.Range("J86").Value = 20004970
$ws.Range("K86").Value = 1733.6666
$ws.Range("L86").Value = 20004970
$ws.Range("M86").Value = -610.6666

$ws.Range("H89").Value = 12503757
$ws.Range("I89").Value = 1733.6666
$ws.Range("J89").Value = 20004970
$ws.Range("K89").Value = 8668.333000000001
$ws.Range("L89").Value = 100024850
$ws.Range("M89").Value = -3052.333000000001

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100:N100").ClearContents()

$ws.Range("H107").Value = 389.95
$ws.Range("I107").Value = 357.8421
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 357.8421
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1562.1579

$ws.Range("H111").Value = 12056.846
$ws.Range("I111").Value = 16548.125
$ws.Range("J111").Value = 4870.8
$ws.Range("K111").Value = 49644.375
$ws.Range("L111").Value = 14612.4
$ws.Range("M111").Value = -46577.375
$ws.Range("N111").Value = -20746.4

$ws.Range("H113").Value = 68463.60000000001
$ws.Range("I113").Value = 201421
$ws.Range("J113").Value = 1984.9
$ws.Range("K113").Value = 201421
$ws.Range("L113").Value = 1984.9
$ws.Range("M113").Value = -198167
$ws.Range("N113").Value = -8492.9

$ws.Range("H115").Value = 33586.668
$ws.Range("I115").Value = 33586.668
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 100760.004
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -99193.00399999999

$ws.Range("H116").Value = 4626.3076
$ws.Range("I116").Value = 5926.75
$ws.Range("J116").Value = 2545.6
$ws.Range("K116").Value = 5926.75
$ws.Range("L116").Value = 2545.6
$ws.Range("M116").Value = -2484.75
$ws.Range("N116").Value = -9429.6

$ws.Range("H137").Value = 1944.6666
$ws.Range("I137").Value = 1592.9286
$ws.Range("J137").Value = 3175.75
$ws.Range("K137").Value = 4778.7858
$ws.Range("L137").Value = 9527.25
$ws.Range("M137").Value = -2228.7858
$ws.Range("N137").Value = -14627.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 60571
$ws.Range("I45").Value = 91977.27
$ws.Range("J45").Value = 2992.8333
$ws.Range("K45").Value = 91977.27
$ws.Range("L45").Value = 2992.8333
$ws.Range("M45").Value = -91600.27
$ws.Range("N45").Value = -3746.8333

$ws.Range("H61").Value = 2408.25
$ws.Range("I61").Value = 1724.75
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 1724.75
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -1512.75
$ws.Range("N61").Value = -3174

$ws.Range("H102").Value = 203252
$ws.Range("I102").Value = 502490
$ws.Range("J102").Value = 3760
$ws.Range("K102").Value = 502490
$ws.Range("L102").Value = 3760
$ws.Range("M102").Value = -500868
$ws.Range("N102").Value = -7004

$ws.Range("H136").Value = 2408.25
$ws.Range("I136").Value = 1724.75
$ws.Range("J136").Value = 2750
$ws.Range("K136").Value = 5174.25
$ws.Range("L136").Value = 8250
$ws.Range("M136").Value = -2624.25
$ws.Range("N136").Value = -13350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3179.0352
$ws.Range("I134").Value = 2972.4048
$ws.Range("J134").Value = 3757.6
$ws.Range("K134").Value = 8917.214399999999
$ws.Range("L134").Value = 11272.8
$ws.Range("M134").Value = -6382.214399999999
$ws.Range("N134").Value = -16342.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37717.15
$ws.Range("I31").Value = 1171.2
$ws.Range("J31").Value = 49899.133
$ws.Range("K31").Value = 1171.2
$ws.Range("L31").Value = 49899.133
$ws.Range("M31").Value = -876.2

$ws.Range("H34").Value = 37717.15
$ws.Range("I34").Value = 1171.2
$ws.Range("J34").Value = 49899.133
$ws.Range("K34").Value = 1171.2
$ws.Range("L34").Value = 49899.133
$ws.Range("M34").Value = -969.2

$ws.Range("H99").Value = 13621.272
$ws.Range("I99").Value = 3132
$ws.Range("J99").Value = 26208.4
$ws.Range("K99").Value = 3132
$ws.Range("L99").Value = 26208.4
$ws.Range("M99").Value = -1634
$ws.Range("N99").Value = -29204.4

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105:N105").ClearContents()

$ws.Range("H126").Value = 13621.272
$ws.Range("I126").Value = 3132
$ws.Range("J126").Value = 26208.4
$ws.Range("K126").Value = 9396
$ws.Range("L126").Value = 78625.20000000001
$ws.Range("M126").Value = -6926
$ws.Range("N126").Value = -83565.20000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 134.5
$ws.Range("I47").Value = 118.125
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 354.375
$ws.Range("L47").Value = 600
$ws.Range("M47").Value = 76.625
$ws.Range("N47").Value = -1462

$ws.Range("H57").Value = 4500
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 4500
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 13500
$ws.Range("N57").Value = -14618
$ws.Range("M57").ClearContents()

$ws.Range("H131").Value = 873.28
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 912.17584
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 2736.52752
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -12816.52752

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3615.6924
$ws.Range("I126").Value = 3763
$ws.Range("J126").Value = 3380
$ws.Range("K126").Value = 11289
$ws.Range("L126").Value = 10140
$ws.Range("M126").Value = -8819
$ws.Range("N126").Value = -15080

$ws.Range("H132").Value = 2160.4348
$ws.Range("I132").Value = 1502.3235
$ws.Range("J132").Value = 4025.0833
$ws.Range("K132").Value = 4506.970499999999
$ws.Range("L132").Value = 12075.2499
$ws.Range("M132").Value = -1976.970499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4353.75
$ws.Range("I7").Value = 1766.6666
$ws.Range("J7").Value = 5906
$ws.Range("K7").Value = 1766.6666
$ws.Range("L7").Value = 5906
$ws.Range("M7").Value = -1654.6666
$ws.Range("N7").Value = -6130

$ws.Range("H126").Value = 4353.75
$ws.Range("I126").Value = 1766.6666
$ws.Range("J126").Value = 5906
$ws.Range("K126").Value = 5299.9998
$ws.Range("L126").Value = 17718
$ws.Range("M126").Value = -2829.9998
$ws.Range("N126").Value = -22658

$ws.Range("H132").Value = 3680.4443
$ws.Range("I132").Value = 3599.375
$ws.Range("J132").Value = 4329
$ws.Range("K132").Value = 10798.125
$ws.Range("L132").Value = 12987
$ws.Range("M132").Value = -8268.125
$ws.Range("N132").Value = -18047

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 100427.5
$ws.Range("I107").Value = 424.375
$ws.Range("J107").Value = 500440
$ws.Range("K107").Value = 1273.125
$ws.Range("L107").Value = 1501320
$ws.Range("M107").Value = 646.875
$ws.Range("N107").Value = -1505160

$ws.Range("H113").Value = 790.5
$ws.Range("I113").Value = 556.5
$ws.Range("J113").Value = 1024.5
$ws.Range("K113").Value = 1669.5
$ws.Range("L113").Value = 3073.5
$ws.Range("M113").Value = 500.5
$ws.Range("N113").Value = -7413.5
